$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet Sheet1 -> sheet1
$ws.Name = "sheet1"

# Give column A a width (closest achievable to the authored 7.85546875 char width)
$ws.Columns.Item(1).ColumnWidth = 7.0

# Add "Numero" primary-key column (column A) for the existing 5 rows
$ws.Range("A1").Value = "Numero"
$ws.Range("A2").Value = 123
$ws.Range("A3").Value = 456
$ws.Range("A4").Value = 789
$ws.Range("A5").Value = 21

# New row 6
$ws.Range("A6").Value = 21
$ws.Range("B6").Value = "dari"
$ws.Range("C6").Value = "hala"
$ws.Range("D6").Value = 123456789
$ws.Range("E6").Value = "haja@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:haja@gmail.com")
$ws.Range("E6").Style = "Hyperlink"

# New row 7
$ws.Range("A7").Value = 123
$ws.Range("B7").Value = "gg"
$ws.Range("C7").Value = "dfd"
$ws.Range("D7").Value = 1233
$ws.Range("E7").Value = "gg@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:gg@gmail.com")
$ws.Range("E7").Style = "Hyperlink"

# Match the selection left by the author on the last populated cell
[void]$ws.Range("E7").Select()
